$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$staging = $ws.Range("ZZ1")
$staging.NumberFormat = "@"

function Set-TextValue($a1ref, $val) {
    $staging.Value = $val
    $staging.Copy()
    $ws.Range($a1ref).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
}

Set-TextValue 'D2' '28.525.12'

Set-TextValue 'D3' '1.867.02'
$ws.Range('E3').Value = '  +0.03%  '

$ws.Range('E4').Value = '  -0.10%  '

Set-TextValue 'D5' '324.81'
$ws.Range('E5').Value = '  -0.58%  '

Set-TextValue 'D6' '1.006'
$ws.Range('E6').Value = '  +0.10%  '

Set-TextValue 'D7' '0.4543'
$ws.Range('E7').Value = '  -1.83%  '

$ws.Range('E8').Value = '  -1.92%  '

Set-TextValue 'D9' '0.07823'
$ws.Range('E9').Value = '  -0.93%  '

Set-TextValue 'D10' '0.9898'
$ws.Range('E10').Value = '  +2.03%  '

Set-TextValue 'D11' '21.50'
$ws.Range('E11').Value = '  -3.55%  '

Set-TextValue 'D12' '1.903.86'
$ws.Range('E12').Value = '  +4.57%  '

Set-TextValue 'D13' '6.924'
$ws.Range('E13').Value = '  -0.12%  '

Set-TextValue 'D14' '5.637'

Set-TextValue 'D15' '0.06932'
$ws.Range('E15').Value = '  +0.06%  '

$ws.Range('B16').Value = 'BinanceUSD'
$ws.Range('C16').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue 'D16' '1.007'
$ws.Range('E16').Value = '  +0.08%  '

$ws.Range('B17').Value = 'Litecoin'
$ws.Range('C17').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 'D17' '86.45'
$ws.Range('E17').Value = '  -2.36%  '

Set-TextValue 'D18' '0.000009942'
$ws.Range('E18').Value = '  -0.98%  '

$ws.Range('E19').Value = '  -1.51%  '

Set-TextValue 'D20' '1.006'
$ws.Range('E20').Value = '  +0.17%  '

Set-TextValue 'D21' '28.531.91'
$ws.Range('E21').Value = '  -0.45%  '

Set-TextValue 'D22' '5.253'
$ws.Range('E22').Value = '  -1.23%  '

Set-TextValue 'D23' '10.89'
$ws.Range('E23').Value = '  -1.65%  '

$ws.Range('B24').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C24').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue 'D24' '2.137.41'
$ws.Range('E24').Value = '  +1.13%  '

$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D25' '2.088'
$ws.Range('E25').Value = '  -1.66%  '

Set-TextValue 'D26' '153.50'
$ws.Range('E26').Value = '  -1.11%  '

Set-TextValue 'D27' '19.06'
$ws.Range('E27').Value = '  -1.29%  '

Set-TextValue 'D28' '5.686'
$ws.Range('E28').Value = '  -1.54%  '

Set-TextValue 'D29' '117.26'
$ws.Range('E29').Value = '  -1.57%  '

Set-TextValue 'D30' '1.879'
$ws.Range('E30').Value = '  -5.73%  '

Set-TextValue 'D31' '0.09270'
$ws.Range('E31').Value = '  -0.81%  '

Set-TextValue 'D32' '0.9047'
$ws.Range('E32').Value = '  -3.49%  '

Set-TextValue 'D33' '5.273'
$ws.Range('E33').Value = '  -0.88%  '

$ws.Range('E34').Value = '  -1.40%  '

Set-TextValue 'D35' '3.265'
$ws.Range('E35').Value = '  -2.49%  '

Set-TextValue 'D36' '0.05660'
$ws.Range('E36').Value = '  -2.72%  '

$ws.Range('E37').Value = '  -0.66%  '

Set-TextValue 'D38' '0.02040'
$ws.Range('E38').Value = '  -3.51%  '

Set-TextValue 'D39' '7.615'
$ws.Range('E39').Value = '  -3.37%  '

Set-TextValue 'D40' '0.5566'
$ws.Range('E40').Value = '  -1.54%  '

Set-TextValue 'D41' '0.1763'
$ws.Range('E41').Value = '  -0.83%  '

Set-TextValue 'D42' '9.648'
$ws.Range('E42').Value = '  -2.86%  '

Set-TextValue 'D43' '0.07134'
$ws.Range('E43').Value = '  -1.63%  '

$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D44' '11.56'
$ws.Range('E44').Value = '  -1.40%  '

$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue 'D45' '0.5240'
$ws.Range('E45').Value = '  -1.41%  '

Set-TextValue 'D46' '2.132'
$ws.Range('E46').Value = '  -3.39%  '

Set-TextValue 'D47' '1.118'
$ws.Range('E47').Value = '  -1.98%  '

Set-TextValue 'D48' '1.804'
$ws.Range('E48').Value = '  -2.31%  '

Set-TextValue 'D49' '111.64'
$ws.Range('E49').Value = '  -1.69%  '

Set-TextValue 'D50' '2.438'
$ws.Range('E50').Value = '  +3.73%  '

Set-TextValue 'D51' '1.005'
$ws.Range('E51').Value = '  +0.03%  '

$staging.Clear()